$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the zoom level of the active window/sheet view
$excel.ActiveWindow.Zoom = 85

# New data values for rows 2-11 (columns A:D)
$data = @(
    @(1, 19.847000000000001, 6, 14),
    @(2, 14.95, 6, 14),
    @(3, 12.105, 4, 13),
    @(4, 16.033999999999999, 8, 14),
    @(5, 17.843, 8, 16),
    @(6, 22.361999999999998, 5, 14),
    @(7, 10.693, 6, 14),
    @(8, 19.23, 6, 15),
    @(9, 11.377000000000001, 7, 14),
    @(10, 10.561, 7, 14)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $r = $i + 2
    $row = $data[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
}

# Remove now-unused rows 12-17 (previously rows 12-17 held entries 11-16)
$ws.Range("A12:D17").ClearContents()
